$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.754.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.345.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.04'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.634'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.613'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.46'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0927'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.49'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.997'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.83'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.704.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.353.99'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.782.06'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.60'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.50'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.55'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '267.94'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.32'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.28'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.44'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.32'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.93'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0904'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.08'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.121'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +15.92%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.64'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0358'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.79'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.65'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -8.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.26'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +17.89%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.36'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.80%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.33'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.03'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '79.09'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +16.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.50'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.06'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.14%  '
